# "add tabel format baru" — bump the table numbering (Tabel 4.2.3/4/5 -> 4.2.5/6/7)
# and the reporting year (2020 -> 2021) on the "Bab 4" sheet's four side-by-side
# mini-tables (columns A/H/P/W), matching the new "013 Polinggona" blanko layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (columns A:G) header stays "Tabel 4.2.2" — unchanged ---

# Title for table 1 (col B, merged B1:G1): year 2020 -> 2021
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Polinggona. 2021"

# --- Table 2 (columns H:M) ---
# Own label: "Tabel 4.2.3" -> "Tabel 4.2.5" (plain text, no rich-text run split)
$ws.Range("H1").Value = "Tabel 4.2.5"

# Title for table 2 (col I, merged I1:M1): year 2020 -> 2021
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Polinggona, 2021"

# --- Table 3 (columns P:S) ---
# Own label is rich text: "Tabel" (default run) + " 4.2.4." (sz9/Calibri run) -> " 4.2.6."
$ws.Range("P1").Value = "Tabel 4.2.6."
$p1Suffix = $ws.Range("P1").Characters(6, 7)
$p1Suffix.Font.Name = "Calibri"
$p1Suffix.Font.Size = 9
$p1Suffix.Font.Underline = $false
$p1Suffix.Font.ColorIndex = 1

# Title for table 3 (col Q, merged Q1:S1): year 2020 -> 2021
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Polinggona, 2021"

# --- Table 4 (columns W:Z) ---
# Own label is rich text: "Tabel" (default run) + " 4.2.5." (sz9/Calibri run) -> " 4.2.7."
$ws.Range("W1").Value = "Tabel 4.2.7."
$w1Suffix = $ws.Range("W1").Characters(6, 7)
$w1Suffix.Font.Name = "Calibri"
$w1Suffix.Font.Size = 9
$w1Suffix.Font.Underline = $false
$w1Suffix.Font.ColorIndex = 1

# Title for table 4 (col X, merged X1:Z1): year 2020 -> 2021
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Polinggona, 2021"
